$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.066.70'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '2.647.71'
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.75'
$ws.Range("E5").Value = '  +4.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.71'
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("E10").Value = '  +4.80%  '
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '3.112.42'
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").Value = '61.055.51'
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.08'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("E16").Value = '  +2.63%  '
$ws.Range("D17").Value = '2.661.11'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '354.88'
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("E21").Value = '  +1.51%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.73'
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  +1.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '0.0₃0862'
$ws.Range("E27").Value = '  +2.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.39'
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.18'
$ws.Range("E30").Value = '  +7.40%  '
$ws.Range("E31").Value = '  +4.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.56'
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.04'
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("E34").Value = '  +4.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.918'
$ws.Range("E36").Value = '  +8.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.899'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '309.04'
$ws.Range("E38").Value = '  +4.72%  '
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("E40").Value = '  +1.85%  '
$ws.Range("E41").Value = '  +3.78%  '
$ws.Range("E42").Value = '  +1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0563'
$ws.Range("E43").Value = '  +1.41%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  +3.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.93'
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("E47").Value = '  +2.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.29'
$ws.Range("E48").Value = '  +8.08%  '
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("D50").Value = '1.992.36'
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("E51").Value = '  +2.77%  '
